$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column D (EMPLOYEE_ID),
# shifting EMPLOYEE_ID, MANAGER_ID, DESIGNATION, STAFF, EMPLOYEE_STATUS,
# PROCESS and DEPARTMENT one column to the right (D->E, E->F, F->G, G->H,
# H->I, I->J, J->K).
$ws.Columns.Item(4).Insert()

# Give the newly inserted column a header; the data cells beneath it are
# left blank (unnamed index column).
$ws.Cells.Item(1,4).Value = "Unnamed: 0"

# Correct the MANAGER_ID (column F), EMPLOYEE_STATUS (column I), PROCESS
# (column J) and DEPARTMENT (column K) values for the affected rows.
$ws.Cells.Item(2,6).Value = ""
$ws.Cells.Item(2,9).Value = "INACTIVE"
$ws.Cells.Item(3,6).Value = "O254"
$ws.Cells.Item(3,10).Value = "MUTHOOT"
$ws.Cells.Item(4,6).Value = "O50"
$ws.Cells.Item(5,6).Value = "O50"
$ws.Cells.Item(6,6).Value = "O254"
$ws.Cells.Item(8,6).Value = "O72"
$ws.Cells.Item(9,6).Value = ""
$ws.Cells.Item(10,6).Value = "O72"
$ws.Cells.Item(10,10).Value = "L&T"
$ws.Cells.Item(11,10).Value = "IDFC"
$ws.Cells.Item(11,11).Value = "TW"
$ws.Cells.Item(12,6).Value = "O254"
$ws.Cells.Item(12,11).Value = "TW"
$ws.Cells.Item(13,9).Value = "INACTIVE"
$ws.Cells.Item(14,6).Value = "O239"
$ws.Cells.Item(16,6).Value = "O50"
$ws.Cells.Item(18,6).Value = "O72"
$ws.Cells.Item(18,10).Value = "MUTHOOT"
$ws.Cells.Item(19,6).Value = "O72"
$ws.Cells.Item(21,11).Value = "RECOVERY"
$ws.Cells.Item(23,6).Value = "O72"
$ws.Cells.Item(24,6).Value = "O50"
$ws.Cells.Item(25,6).Value = "O72"
$ws.Cells.Item(26,6).Value = ""
$ws.Cells.Item(27,6).Value = "O50"
